# PLANO_TRABALHO_CJF_PROSA.docx — applies the edits described in the commit diff.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "ingressou em licença médica imediata" -> "ingressou em licença médica"
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("ingressou em licença médica imediata", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ingressou em licença médica", 2)
Write-Output "1) licença médica: $ok"

# ---------------------------------------------------------------------------
# 2) Bold "primeiro lugar " inside the Meta 1 sentence, splitting it into the
#    4 runs shown in the diff ( ". A conquista do " / "primeiro lugar " (bold) /
#    "no cumprimento ... Justiça" / ", com índice ... implementada." )
# ---------------------------------------------------------------------------
$rngBold = $d.Content
$ok = $rngBold.Find.Execute("primeiro lugar ", $true, $false, $false, $false, $false,
                             $true, 1, $false, "", 0)
Write-Output "2) found 'primeiro lugar ': $ok"
$rngBold.Font.Bold = $true

# Touching (and reverting) the following span forces Word to keep it as its
# own run, split away from the still-unbolded remainder of the sentence —
# matching the 4-run layout in the target XML while leaving its formatting
# (just pt-BR lang) untouched.
$rngSplit = $d.Content
$ok = $rngSplit.Find.Execute("no cumprimento da Meta 1 do Conselho Nacional de Justiça", $true, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)
Write-Output "3) found trailing span: $ok"
$rngSplit.Font.Bold = $true
$rngSplit.Font.Bold = $false

# ---------------------------------------------------------------------------
# 3) Highlight (yellow) the IAD paragraph
# ---------------------------------------------------------------------------
$rngIAD = $d.Content
$ok = $rngIAD.Find.Execute("O Índice de Atendimento à Demanda (IAD), que mensura a relação entre processos baixados e casos novos distribuídos, evidencia a reversão completa do quadro de insuficiência produtiva. A evolução de 41,26% para 113,92% representa incremento de 176 pontos percentuais, indicando que a unidade não apenas equiparou sua capacidade de processamento à demanda de distribuição, mas passou a operar em ritmo superior ao necessário para estabilização do acervo, permitindo a redução ativa do estoque acumulado.",
                             $true, $false, $false, $false, $false,
                             $true, 1, $false, "", 0)
Write-Output "4) found IAD paragraph: $ok"
$rngIAD.HighlightColorIndex = 7

# ---------------------------------------------------------------------------
# 4) Highlight (yellow) the TCL paragraph
# ---------------------------------------------------------------------------
$rngTCL = $d.Content
$ok = $rngTCL.Find.Execute("Paralelamente, a Taxa de Congestionamento Líquida (TCL) — indicador que mede o percentual de processos que permaneceram em tramitação sem solução definitiva ao longo do exercício — experimentou redução de 51 pontos percentuais, transitando de 69,76% para 33,86%. Esta redução situa a unidade em patamar inferior à média nacional das Turmas Recursais (que se mantém historicamente acima de 50%), demonstrando não apenas a recuperação, mas a excelência do desempenho atual quando comparado aos padrões nacionais de produtividade judiciária.",
                             $true, $false, $false, $false, $false,
                             $true, 1, $false, "", 0)
Write-Output "5) found TCL paragraph: $ok"
$rngTCL.HighlightColorIndex = 7

# ---------------------------------------------------------------------------
# 5) "documentação fotográfica ... verificação visual" -> drop "fotográfica"/"visual"
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("A documentação fotográfica dos painéis oficiais do CNJ, anexa a este plano, permite a verificação visual da trajetória ascendente dos indicadores ao longo dos meses de 2025, evidenciando a consistência e sustentabilidade dos resultados alcançados.",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "A documentação dos painéis oficiais do CNJ, anexa a este plano, permite a verificação da trajetória ascendente dos indicadores ao longo dos meses de 2025, evidenciando a consistência e sustentabilidade dos resultados alcançados.",
                         2)
Write-Output "6) documentação fotográfica: $ok"

# ---------------------------------------------------------------------------
# 6) "O objetivo estratégico para" -> "O objetivo para"
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("O objetivo estratégico para o exercício de 2026", $true, $false, $false, $false, $false,
                         $true, 1, $false, "O objetivo para o exercício de 2026", 2)
Write-Output "7) objetivo estratégico: $ok"

# ---------------------------------------------------------------------------
# 7) "recurso humano central" -> "recurso central"
#    (the accompanying shift of "demonstrada e permite a condução simultânea
#    do " across the page-break run is a pagination artefact only — the
#    rendered/concatenated paragraph text is identical either way, so it is
#    intentionally left alone to avoid corrupting the lastRenderedPageBreak
#    marker on the following run.)
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("recurso humano central", $true, $false, $false, $false, $false,
                         $true, 1, $false, "recurso central", 2)
Write-Output "8) recurso humano central: $ok"

# ---------------------------------------------------------------------------
# 8) Table cell "18,0 meses" -> "18 meses"
# ---------------------------------------------------------------------------
$rng = $d.Content
$ok = $rng.Find.Execute("18,0 meses", $true, $false, $false, $false, $false,
                         $true, 1, $false, "18 meses", 2)
Write-Output "9) 18,0 meses: $ok"

Write-Output "all edits applied"
